$wb = $excel.ActiveWorkbook

# --- "neumonia" sheet (sheet1): fill in week-14 row with the 2015-2020 values ---
$ws1 = $wb.Worksheets.Item("neumonia")
$ws1.Range("B15").Value = 2859
$ws1.Range("C15").Value = 3386
$ws1.Range("D15").Value = 2995
$ws1.Range("E15").Value = 2128
$ws1.Range("F15").Value = 2945
$ws1.Range("G15").Value = 2181

# --- "ira" sheet (sheet2): fill in week-14 row with the 2015-2020 values ---
$ws2 = $wb.Worksheets.Item("ira")
$ws2.Range("B15").Value = 358028
$ws2.Range("C15").Value = 422634
$ws2.Range("D15").Value = 576005
$ws2.Range("E15").Value = 396706
$ws2.Range("F15").Value = 521737
$ws2.Range("G15").Value = 483034

# --- restore the selections left behind on each sheet ---
$ws2.Range("E22").Select()

# Select on "neumonia" last so it stays the active/visible tab (tabSelected="1")
$ws1.Activate()
$ws1.Range("D15").Select()
